# [FIX] assign multiples coureurs to etape
#
# Rows 27-37 ("liste des etapes profil equipe") get actual "Temps passé"
# (E) values recorded now that multiple coureurs can be assigned to an
# etape, and "Reste à faire" (F) drops back to 0 for each (G recalculates
# via the existing shared formula =(100-F)).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new "Temps passé" (E) value. F is reset to 0 for each of these rows,
# and G (=100-F) will recompute to 100 automatically.
$updates = @{
    27 = 15
    28 = 10
    29 = 5
    30 = 10
    31 = 5
    32 = 10
    33 = 10
    34 = 10
    35 = 10
    36 = 10
    37 = 10
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 5).Value = $updates[$row]
    $ws.Cells.Item($row, 6).Value = 0
}

# Update the saved view state: scroll back to the top (drop topLeftCell)
# and move the active selection to B44.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("B44").Select()
